$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2 = @{ G = 0.1421496666666667; H = 0.426449; I = 0.1211014306728536; J = 0.1211014306728536;
           M = 0.01569233333333333; N = 0.047077; O = 0.03693539111407157; P = 0.03693539111407157;
           Q = 0.002230659952555556; R = 0.020075939573; S = 0.00447292870637547; T = 0.00447292870637547 }
    3 = @{ G = 0.1421496666666667; H = 0.426449; I = 0.1211014306728536; J = 0.1211014306728536;
           N = 0.9690430000000001; O = 0.7602859615386125; P = 0.7602859615386125;
           Q = 0.0459163798118889; R = 0.4132474183070001; S = 0.09207171766281209; T = 0.09207171766281209 }
    4 = @{ G = 0.1421496666666667; H = 0.426449; I = 0.1211014306728536; J = 0.1211014306728536;
           M = 0.08615233333333333; N = 0.258457; O = 0.202778647347316; P = 0.202778647347316;
           Q = 0.01224652546588889; R = 0.110218729193; S = 0.02455678430366601; T = 0.02455678430366601 }
    5 = @{ I = 0.8788985693271465; J = 0.8788985693271465;
           M = 0.01569233333333333; N = 0.047077; O = 0.03693539111407157; P = 0.03693539111407157;
           Q = 0.01618910552966667; R = 0.145701949767; S = 0.0324624624076961; T = 0.0324624624076961 }
    6 = @{ I = 0.8788985693271465; J = 0.8788985693271465;
           N = 0.9690430000000001; O = 0.7602859615386125; P = 0.7602859615386125;
           Q = 0.3332399980836667; S = 0.6682142438758004; T = 0.6682142438758004 }
    7 = @{ I = 0.8788985693271465; J = 0.8788985693271465;
           M = 0.08615233333333333; N = 0.258457; O = 0.202778647347316; P = 0.202778647347316;
           Q = 0.08887965774966666; R = 0.799916919747; S = 0.17822186304365; T = 0.17822186304365 }
}

foreach ($row in $newValues.Keys) {
    foreach ($col in $newValues[$row].Keys) {
        $ws.Range("$col$row").Value = $newValues[$row][$col]
    }
}
